$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet so it becomes "Sheet4"
# and the final, active tab (matches activeTab going from 2 -> 3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet4"

# Header row (labels repeated across top and left to form a symmetric matrix)
# Written before the title so the shared-string table order is
# Area, RL, WL, RE, WE, LP, Bigtable (matches target).
$ws.Range("C3").Value = "Area"
$ws.Range("D3").Value = "RL"
$ws.Range("E3").Value = "WL"
$ws.Range("F3").Value = "RE"
$ws.Range("G3").Value = "WE"
$ws.Range("H3").Value = "LP"

# Data rows
$ws.Range("B4").Value = "Area"
$ws.Range("C4").Value = 3.06
$ws.Range("D4").Value = 10.7
$ws.Range("E4").Value = 16.4
$ws.Range("F4").Value = 5.66
$ws.Range("G4").Value = 6.22
$ws.Range("H4").Value = 3.63

$ws.Range("B5").Value = "RL"
$ws.Range("C5").Value = 21.8
$ws.Range("D5").Value = 3.7
$ws.Range("E5").Value = 4.92
$ws.Range("F5").Value = 9.12
$ws.Range("G5").Value = 9.57
$ws.Range("H5").Value = 9.91

$ws.Range("B6").Value = "WL"
$ws.Range("C6").Value = 18.6
$ws.Range("D6").Value = 13.9
$ws.Range("E6").Value = 4.01
$ws.Range("F6").Value = 15.9
$ws.Range("G6").Value = 11.3
$ws.Range("H6").Value = 18.1

$ws.Range("B7").Value = "RE"
$ws.Range("C7").Value = 0.276
$ws.Range("D7").Value = 0.225
$ws.Range("E7").Value = 0.316
$ws.Range("F7").Value = 0.105
$ws.Range("G7").Value = 0.139
$ws.Range("H7").Value = 0.279

$ws.Range("B8").Value = "WE"
$ws.Range("C8").Value = 0.293
$ws.Range("D8").Value = 0.322
$ws.Range("E8").Value = 0.309
$ws.Range("F8").Value = 0.193
$ws.Range("G8").Value = 0.131
$ws.Range("H8").Value = 0.281

$ws.Range("B9").Value = "LP"
$ws.Range("C9").Value = 1.01
$ws.Range("D9").Value = 3.53
$ws.Range("E9").Value = 4.98
$ws.Range("F9").Value = 1.85
$ws.Range("G9").Value = 1.92
$ws.Range("H9").Value = 0.78

# Title (written last so "Bigtable" lands at the end of the shared-string table).
$ws.Range("B2").Value = "Bigtable"

# Selection matches the commit: active cell B2 on the new sheet.
$ws.Range("B2").Select()
